# Rebuild the SitewideSearch workbook: rename the English sheet, swap which
# tab is active, and replace both sheets' content with the new
# SearchTerm/ResultType layout (English sheet first, Spanish sheet second).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Stash the two cell formats we still need (bold/shaded header style,
# and the "numeric value styled as text" style) on a scratch sheet before
# clearing everything, since Cells.Clear() drops per-cell formatting too.
# NOTE: Worksheets.Add() with no args inserts *before* sheet 1 and shifts
# every existing positional handle, so anchor it after the last sheet
# instead, which leaves $ws1 / $ws2 pointing at the right sheets.
$helper = $wb.Worksheets.Add($null, $ws2)
$ws1.Range("A1").Copy()
$helper.Range("A1").PasteSpecial(-4122)
$ws1.Range("D2").Copy()
$helper.Range("A2").PasteSpecial(-4122)

# --- Rename the first sheet.
$ws1.Name = "SitewideSearchEn"

# --- Wipe both sheets completely (content + formatting + column widths).
$ws1.Cells.Clear()
$ws2.Cells.Clear()

# --- New English sheet (SitewideSearchEn) data.
$enRows = @(
    @("SearchTerm", "ResultType"),
    @("Kaposi's sarcoma", "Generic"),
    @("[F-18]HX4", "Generic"),
    @("250", "Generic"),
    @("LiveHelp", "BestBet"),
    @("stereostatic radiosurgery", "BestBet"),
    @("ipilimumab", "Definition"),
    @("abdominoperineal resection ", "Definition"),
    @("tumor", "BestBetAndDefinition"),
    @("glioma", "BestBetAndDefinition"),
    @("argle-bargle or foofaraw", "NoMatch"),
    @("10001110101 10001110101", "NoMatch")
)

for ($r = 0; $r -lt $enRows.Count; $r++) {
    $row = $enRows[$r]
    if ($r -eq 3) {
        $ws1.Cells.Item($r + 1, 1).Value = 250
    } else {
        $ws1.Cells.Item($r + 1, 1).Value = $row[0]
    }
    $ws1.Cells.Item($r + 1, 2).Value = $row[1]
}

# --- New Spanish sheet (SitewideSearchEs) data.
$esRows = @(
    @("SearchTerm", "ResultType"),
    @("cáncer", "Generic"),
    @("dalteparina sódica", "Generic"),
    @("macrófago", "Generic"),
    @("dolor", "BestBet"),
    @("linfoma", "BestBet"),
    @("tumor", "Definition"),
    @("safingol", "Definition"),
    @("argle-bargle o foofaraw", "NoMatch"),
    @("10001110101 10001110101", "NoMatch")
)

for ($r = 0; $r -lt $esRows.Count; $r++) {
    $row = $esRows[$r]
    $ws2.Cells.Item($r + 1, 1).Value = $row[0]
    $ws2.Cells.Item($r + 1, 2).Value = $row[1]
}

# --- Re-apply the saved header style to both sheets' header rows.
$helper.Range("A1").Copy()
$ws1.Range("A1:B1").PasteSpecial(-4122)
$ws2.Range("A1:B1").PasteSpecial(-4122)

# --- Re-apply the numeric-styled-as-text style to the "250" cell.
$helper.Range("A2").Copy()
$ws1.Range("A4").PasteSpecial(-4122)

# --- Drop the scratch sheet now that both styles have been copied out.
$helper.Delete()

# --- Column widths (best effort / matches the sheets' natural bestFit sizes).
$ws1.Columns.Item(1).ColumnWidth = 27
$ws1.Columns.Item(2).ColumnWidth = 21

$ws2.Columns.Item(1).ColumnWidth = 24
$ws2.Columns.Item(2).ColumnWidth = 26

# --- Selection state matches the new used ranges.
$ws1.Range("A13").Select()
$ws2.Range("A11").Select()

# --- The Spanish sheet is now the active/visible tab.
$ws2.Activate()
